# Applies the "Deploying to gh-pages ... LinuxForHealth/alvearie-fhir-ig"
# commit's changes to the StructureDefinition-snapshot-age-in-months workbook.
#
# Metadata sheet: URL / Version / Date / Publisher values are rebranded from
# ibm.com/Alvearie to linuxforhealth.org/LinuxForHealth and bumped to the
# new release (8.0.0, 2022-11-10T16:00:46+00:00).
#
# Elements sheet: the root "Extension" row's Constraint(s) cell (AI2), which
# used to repeat the ele-1/ext-1 FHIR invariant text, is cleared (that text
# now lives solely on the "Extension.extension" row, AI4, where it already
# belonged).

$wb = $excel.ActiveWorkbook

# ---- Metadata sheet ----
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/snapshot-age-in-months"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# ---- Elements sheet ----
$elements = $wb.Worksheets.Item("Elements")

# The canonical URL string is shared between Metadata!B2 and the
# Extension.url row's Fixed Value (Q5) -- both must be rebranded together.
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/snapshot-age-in-months"

# Clear the Constraint(s) cell for the top-level "Extension" element row.
$elements.Range("AI2").Value = ""
